$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force every written value to be stored as literal text (matches the
# original inlineStr cells) by prefixing a quote -- exactly like typing
# `26.514.00 into a General-formatted cell in the Excel UI -- then restore
# the cell style to Normal so no stray quotePrefix/number-format style
# sticks around on cells that never had one.
$ws.Range("D2").Value = "'26.514.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.36%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.679.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.71%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.21%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'219.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.04%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.5318"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.42%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.22%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2703"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +3.37%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06415"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.24%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'21.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.61%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.28%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.695.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +3.31%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.517"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.89%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.5601"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.03%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0₅8357"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.49%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'65.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.75%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'26.545.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.51%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -0.05%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.807"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.97%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'193.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.04%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'10.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.81%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.347"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.25%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.24%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.1278"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +5.42%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'139.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.44%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'7.424"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.36%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +2.07%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.444"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +3.51%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.06297"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +6.42%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +2.15%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.615"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +4.72%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.466"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.34%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +2.38%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +2.61%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.6203"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +8.93%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.421"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.33%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.790"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.98%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'VeChain"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'0.01635"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.99%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'FraxShare"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'6.147"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +6.96%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.096.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +6.10%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.8629"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.54%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.0000"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.07%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'100.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.32%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.823.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.91%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'BabyDogeCoin"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'0.0₈111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.18%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Aave"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'58.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.88%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'8.168"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.12%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.01%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.486"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +6.57%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.05197"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.75%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'6.045"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.65%  "
$ws.Range("E51").Style = "Normal"
